$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.958.07'
$ws.Range('E2').Value = '  -0.83%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.915.13'
$ws.Range('E3').Value = '  -0.13%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.003'
$ws.Range('E4').Value = '  +0.25%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '320.04'
$ws.Range('E5').Value = '  -2.88%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.003'
$ws.Range('E6').Value = '  +0.35%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5037'
$ws.Range('E7').Value = '  -3.08%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4020'
$ws.Range('E8').Value = '  -1.13%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.08239'
$ws.Range('E9').Value = '  -3.04%  '

# Row 10
$ws.Range('E10').Value = '  -1.54%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '42.03'
$ws.Range('E11').Value = '  -1.90%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '23.53'
$ws.Range('E12').Value = '  +0.60%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.909.16'
$ws.Range('E13').Value = '  -0.43%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.399'
$ws.Range('E14').Value = '  -0.79%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.298'
$ws.Range('E15').Value = '  -1.43%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.004'
$ws.Range('E16').Value = '  +0.33%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '92.17'
$ws.Range('E17').Value = '  -3.34%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001095'
$ws.Range('E18').Value = '  -1.64%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06485'
$ws.Range('E19').Value = '  -3.05%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '18.49'
$ws.Range('E20').Value = '  +0.07%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.003'
$ws.Range('E21').Value = '  +0.37%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.952'
$ws.Range('E22').Value = '  -1.06%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '29.994.12'
$ws.Range('E23').Value = '  -0.76%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.27'
$ws.Range('E24').Value = '  -0.73%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.194'
$ws.Range('E25').Value = '  -1.79%  '

# Row 26
$ws.Range('B26').Value = 'EthereumClassic'
$ws.Range('C26').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '22.06'
$ws.Range('E26').Value = '  +2.87%  '

# Row 27
$ws.Range('B27').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C27').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.131.79'
$ws.Range('E27').Value = '  -0.36%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '161.09'
$ws.Range('E28').Value = '  -0.50%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.335'
$ws.Range('E29').Value = '  -3.51%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '128.69'
$ws.Range('E30').Value = '  -0.17%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.122'
$ws.Range('E31').Value = '  +1.93%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.1041'
$ws.Range('E32').Value = '  -2.50%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.980'
$ws.Range('E33').Value = '  -0.68%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.815'
$ws.Range('E34').Value = '  +4.92%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.02442'
$ws.Range('E35').Value = '  -2.12%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.371'
$ws.Range('E36').Value = '  +3.52%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.06436'
$ws.Range('E37').Value = '  -2.28%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2155'
$ws.Range('E38').Value = '  -2.67%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '8.866'
$ws.Range('E39').Value = '  +0.40%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.190'
$ws.Range('E40').Value = '  -3.19%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.6419'
$ws.Range('E41').Value = '  -1.73%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '11.37'
$ws.Range('E42').Value = '  -4.39%  '

# Row 43
$ws.Range('E43').Value = '  -2.06%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.003'
$ws.Range('E44').Value = '  +0.37%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.23'
$ws.Range('E45').Value = '  -0.46%  '

# Row 46
$ws.Range('B46').Value = 'Decentraland'
$ws.Range('C46').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5990'
$ws.Range('E46').Value = '  -2.67%  '

# Row 47
$ws.Range('B47').Value = 'NEARProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.158'
$ws.Range('E47').Value = '  +3.66%  '

# Row 48
$ws.Range('E48').Value = '  -2.64%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '122.85'
$ws.Range('E49').Value = '  -1.13%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.213'
$ws.Range('E50').Value = '  -2.67%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '78.83'
$ws.Range('E51').Value = '  -0.93%  '
